# Apply updated crypto price/volume figures to worksheet cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''30.737.08'
$ws.Range("D3").Value = '''1.894.72'
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '''246.20'
$ws.Range("E5").Value = '  +2.01%  '
$ws.Range("D6").Value = '''0.9994'
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").Value = '''0.4928'
$ws.Range("E7").Value = '  -1.29%  '
$ws.Range("D8").Value = '''0.2957'
$ws.Range("E8").Value = '  +1.25%  '
$ws.Range("D9").Value = '''0.06809'
$ws.Range("E9").Value = '  +2.99%  '
$ws.Range("D10").Value = '''1.896.58'
$ws.Range("E10").Value = '  +1.00%  '
$ws.Range("E11").Value = '  +3.90%  '
$ws.Range("D12").Value = '''92.29'
$ws.Range("E12").Value = '  +7.07%  '
$ws.Range("D13").Value = '''0.07267'
$ws.Range("E13").Value = '  +0.25%  '
$ws.Range("D14").Value = '''0.6856'
$ws.Range("E14").Value = '  +2.83%  '
$ws.Range("D15").Value = '''5.098'
$ws.Range("E15").Value = '  +4.77%  '
$ws.Range("D16").Value = '''30.712.43'
$ws.Range("E16").Value = '  +2.63%  '
$ws.Range("D17").Value = '''0.000008013'
$ws.Range("E17").Value = '  +1.51%  '
$ws.Range("D18").Value = '''13.37'
$ws.Range("E18").Value = '  +4.79%  '
$ws.Range("D19").Value = '''0.9996'
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("D20").Value = '''2.136.09'
$ws.Range("E20").Value = '  +0.85%  '
$ws.Range("D21").Value = '''1.002'
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("D22").Value = '''4.853'
$ws.Range("E22").Value = '  +2.01%  '
$ws.Range("D23").Value = '''194.17'
$ws.Range("E23").Value = '  +39.86%  '
$ws.Range("D24").Value = '''6.090'
$ws.Range("E24").Value = '  +8.19%  '
$ws.Range("D25").Value = '''9.352'
$ws.Range("E25").Value = '  +3.12%  '
$ws.Range("D26").Value = '''155.31'
$ws.Range("E26").Value = '  +4.35%  '
$ws.Range("D27").Value = '''19.49'
$ws.Range("E27").Value = '  +14.86%  '
$ws.Range("D28").Value = '''1.933'
$ws.Range("E28").Value = '  +1.30%  '
$ws.Range("D29").Value = '''1.395'
$ws.Range("E29").Value = '  +0.79%  '
$ws.Range("D30").Value = '''4.344'
$ws.Range("E30").Value = '  +3.63%  '
$ws.Range("D31").Value = '''0.09026'
$ws.Range("E31").Value = '  +2.63%  '
$ws.Range("D32").Value = '''4.049'
$ws.Range("E32").Value = '  +2.48%  '
$ws.Range("D33").Value = '''0.05194'
$ws.Range("E33").Value = '  +2.80%  '
$ws.Range("D34").Value = '''0.7508'
$ws.Range("E34").Value = '  +5.47%  '
$ws.Range("E35").Value = '  +2.14%  '
$ws.Range("E36").Value = '  +1.30%  '
$ws.Range("D37").Value = '''0.01883'
$ws.Range("E37").Value = '  +7.90%  '
$ws.Range("D38").Value = '''2.669'
$ws.Range("E38").Value = '  -0.98%  '
$ws.Range("E39").Value = '  -0.58%  '
$ws.Range("D40").Value = '''0.9391'
$ws.Range("E40").Value = '  +1.00%  '
$ws.Range("D41").Value = '''0.4451'
$ws.Range("E41").Value = '  +4.42%  '
$ws.Range("D42").Value = '''105.97'
$ws.Range("E42").Value = '  +4.29%  '
$ws.Range("D43").Value = '''5.830'
$ws.Range("E43").Value = '  +0.82%  '
$ws.Range("D44").Value = '''0.9998'
$ws.Range("D45").Value = '''7.737'
$ws.Range("E45").Value = '  +4.05%  '
$ws.Range("D46").Value = '''0.1344'
$ws.Range("E46").Value = '  +7.06%  '
$ws.Range("D47").Value = '''0.05862'
$ws.Range("E47").Value = '  +3.60%  '
$ws.Range("D48").Value = '''8.774'
$ws.Range("E48").Value = '  +7.18%  '
$ws.Range("D49").Value = '''0.3976'
$ws.Range("E49").Value = '  +5.94%  '
$ws.Range("D50").Value = '''33.63'
$ws.Range("E50").Value = '  +4.03%  '
$ws.Range("D51").Value = '''1.412'
$ws.Range("E51").Value = '  +5.99%  '
